$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarterly row (01-07-2021) appended after the existing data (row 76).
# Format the date-label cell as Text first so the "dd-mm-yyyy"-looking
# string is stored as a literal string (matching the other "Serie" column
# cells) instead of being auto-converted into a date serial number, then
# drop back to the sheet's normal (unstyled) cell style.
$dateCell = $ws.Range("A76")
$dateCell.NumberFormat = "@"
$dateCell.Value = "01-07-2021"
$dateCell.Style = "Normal"

$ws.Range("B76").Value = 2196
$ws.Range("C76").Value = 189
$ws.Range("D76").Value = 172
$ws.Range("E76").Value = 2007
$ws.Range("F76").Value = 680
$ws.Range("G76").Value = 1327
